$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.817.56'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '2.729.24'
$ws.Range('E3').Value = '  +3.15%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'603.09"
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('D6').Value = "'168.99"
$ws.Range('E6').Value = '  +5.95%  '
$ws.Range('E8').Value = '  +0.64%  '
$ws.Range('D9').Value = '2.728.15'
$ws.Range('E9').Value = '  +3.19%  '
$ws.Range('D10').Value = "'0.144"
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('D11').Value = "'0.369"
$ws.Range('E11').Value = '  +5.20%  '
$ws.Range('D12').Value = "'5.34"
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('E14').Value = '  +2.67%  '
$ws.Range('D15').Value = '3.232.23'
$ws.Range('E15').Value = '  +3.27%  '
$ws.Range('E16').Value = '  +1.16%  '
$ws.Range('D17').Value = '68.694.25'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '2.726.54'
$ws.Range('E18').Value = '  +3.89%  '
$ws.Range('D19').Value = "'11.89"
$ws.Range('E19').Value = '  +4.50%  '
$ws.Range('D20').Value = "'373.35"
$ws.Range('E20').Value = '  +3.47%  '
$ws.Range('E21').Value = '  +4.27%  '
$ws.Range('D22').Value = "'4.54"
$ws.Range('E22').Value = '  +3.12%  '
$ws.Range('E23').Value = '  +4.26%  '
$ws.Range('E24').Value = '  +3.01%  '
$ws.Range('D25').Value = "'73.50"
$ws.Range('E25').Value = '  -2.21%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').Value = "'10.02"
$ws.Range('E27').Value = '  +1.81%  '
$ws.Range('D28').Value = '2.872.92'
$ws.Range('E28').Value = '  +3.52%  '
$ws.Range('E29').Value = '  +1.89%  '
$ws.Range('D30').Value = "'588.68"
$ws.Range('E30').Value = '  +4.54%  '
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('E32').Value = '  +3.82%  '
$ws.Range('E33').Value = '  +4.68%  '
$ws.Range('E34').Value = '  +5.03%  '
$ws.Range('E35').Value = '  +3.92%  '
$ws.Range('E36').Value = '  +4.40%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').Value = "'162.19"
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('D39').Value = "'19.95"
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('D41').Value = "'1.92"
$ws.Range('E41').Value = '  +3.08%  '
$ws.Range('E42').Value = '  +3.03%  '
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('D44').Value = "'2.65"
$ws.Range('E44').Value = '  +1.16%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = "'41.02"
$ws.Range('E46').Value = '  +1.91%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0310'
$ws.Range('E47').Value = '  -3.89%  '
$ws.Range('D48').Value = "'156.18"
$ws.Range('E48').Value = '  -0.56%  '
$ws.Range('D49').Value = "'3.97"
$ws.Range('E49').Value = '  +5.62%  '
$ws.Range('E50').Value = '  +6.54%  '
$ws.Range('D51').Value = "'0.604"
$ws.Range('E51').Value = '  +6.70%  '

# Restore default (Normal) style on cells that needed a text-prefix
# to avoid Excel auto-converting numeric-looking strings to numbers.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
